# Add the 2022-Q3 quarterly snapshot:
#  - new worksheet "2022-Q3" (fund holdings detail), inserted before "2022-Q2"
#  - new summary row on "总计" for the 2022-Q3 totals

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)      # "总计"
$q2     = $wb.Worksheets.Item(2)      # "2022-Q2" (existing detail sheet)

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right before the "2022-Q2" tab.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Pull the header style (bold + bordered, same look as the "总计" header row)
# across onto the new sheet's header row and the "A" index column.
$totals.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2:A9").PasteSpecial(-4122)

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Fund holdings detail rows. Columns B, D, E, F, G are text in the source
# data (fund code keeps leading zeros, numeric-looking figures are stored
# as plain text, not numbers) -- force text storage via NumberFormat "@"
# and clear it back off afterwards so no residual style sticks around.
$rows = @(
    @(0, "012390", "中欧产业前瞻混合A",                               "16.91", "91.97", "2.65", "0.4481", 10),
    @(1, "012557", "中欧景气前瞻一年持有期混合型证券投资基金A",       "8.43",  "92.09", "2.67", "0.2251", 10),
    @(2, "007132", "长城港股通价值精选多策略混合",                   "0.72",  "80.94", "5.19", "0.0374", 4),
    @(3, "012558", "中欧景气前瞻一年持有期混合型证券投资基金C",       "0.68",  "92.09", "2.67", "0.0182", 10),
    @(4, "161620", "融通核心价值混合（QDII）A",                      "0.55",  "57.96", "3.23", "0.0178", 5),
    @(5, "012391", "中欧产业前瞻混合C",                               "0.52",  "91.97", "2.65", "0.0138", 10),
    @(6, "003279", "融通沪港深智慧生活灵活配置混合",                 "0.10",  "55.71", "5.56", "0.0056", 3),
    @(7, "014127", "融通核心价值混合（QDII）C",                      "0.01",  "57.96", "3.23", "0.0003", 5)
)

$r = 2
foreach ($row in $rows) {
    $textRange = $q3.Range("B$r,D$r:G$r")
    $textRange.NumberFormat = "@"

    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]

    $textRange.ClearFormats()
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. "总计": push the existing (2022-Q2) row down and insert the new
#    2022-Q3 summary row above it.
# ---------------------------------------------------------------------------
$totals.Rows.Item(2).Insert()
$totals.Range("B2:D2").ClearFormats()

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 8
$totals.Range("D2").Value = 0.77

$totals.Range("A3").Value = 1
